$wb = $excel.ActiveWorkbook

# zh-cn sheet: update the handoff/handback datetimes for the first file row
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-20 10:37:12"
$wsZh.Range("H2").Value = "2016-03-20 10:37:32"

# de-de sheet: update the handoff/handback datetimes for the first file row
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-20 10:37:16"
$wsDe.Range("H2").Value = "2016-03-20 10:37:38"
